# Update "想去人数" (column F) figures across the three affected sheets.
# Sheet "演出" is not touched by this change.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (Exhibitions) sheet ---
$wsExhibit.Range("F4").Value  = 1864
$wsExhibit.Range("F5").Value  = 3238
$wsExhibit.Range("F7").Value  = 4723
$wsExhibit.Range("F9").Value  = 269
$wsExhibit.Range("F11").Value = 613
$wsExhibit.Range("F13").Value = 15
$wsExhibit.Range("F15").Value = 642
$wsExhibit.Range("F16").Value = 281
$wsExhibit.Range("F18").Value = 87
$wsExhibit.Range("F19").Value = 138
$wsExhibit.Range("F20").Value = 331
$wsExhibit.Range("F21").Value = 4685
$wsExhibit.Range("F24").Value = 9
$wsExhibit.Range("F25").Value = 5834
$wsExhibit.Range("F27").Value = 1177
$wsExhibit.Range("F28").Value = 235
$wsExhibit.Range("F29").Value = 651
$wsExhibit.Range("F30").Value = 4408
$wsExhibit.Range("F31").Value = 2
$wsExhibit.Range("F32").Value = 74
$wsExhibit.Range("F33").Value = 116
$wsExhibit.Range("F34").Value = 807
$wsExhibit.Range("F35").Value = 57
$wsExhibit.Range("F36").Value = 734
$wsExhibit.Range("F37").Value = 744

# --- 本地生活 (Local life) sheet ---
$wsLocal.Range("F3").Value = 1075
$wsLocal.Range("F4").Value = 33

# --- 全部类型 (All types) sheet ---
$wsAll.Range("F4").Value  = 1075
$wsAll.Range("F5").Value  = 33
$wsAll.Range("F7").Value  = 1864
$wsAll.Range("F9").Value  = 3238
$wsAll.Range("F11").Value = 4723
$wsAll.Range("F13").Value = 269
$wsAll.Range("F15").Value = 613
$wsAll.Range("F17").Value = 15
$wsAll.Range("F19").Value = 642
$wsAll.Range("F20").Value = 281
$wsAll.Range("F23").Value = 87
$wsAll.Range("F24").Value = 138
$wsAll.Range("F25").Value = 331
$wsAll.Range("F26").Value = 4685
$wsAll.Range("F29").Value = 9
$wsAll.Range("F30").Value = 5834
$wsAll.Range("F32").Value = 1177
$wsAll.Range("F33").Value = 235
$wsAll.Range("F34").Value = 651
$wsAll.Range("F35").Value = 4408
$wsAll.Range("F36").Value = 2
$wsAll.Range("F38").Value = 74
$wsAll.Range("F39").Value = 116
$wsAll.Range("F40").Value = 807
$wsAll.Range("F41").Value = 57
$wsAll.Range("F42").Value = 734
$wsAll.Range("F43").Value = 744
